$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.129.58"
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("D3").Value = "1.856.07"
$ws.Range("E3").Value = "  -2.99%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value2 = "'233.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.54%  "
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value2 = "'0.4697"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.96%  "
$ws.Range("D8").Value2 = "'0.2806"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.39%  "
$ws.Range("D9").Value2 = "'0.06537"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.47%  "
$ws.Range("D10").Value2 = "'19.90"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.28%  "
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D12").Value2 = "'96.68"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.02%  "
$ws.Range("D13").Value = "1.853.44"
$ws.Range("E13").Value = "  -3.29%  "
$ws.Range("D14").Value2 = "'5.083"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.35%  "
$ws.Range("D15").Value2 = "'0.6653"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value2 = "'281.58"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.27%  "
$ws.Range("D17").Value = "30.169.73"
$ws.Range("E17").Value = "  -1.27%  "
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").Value2 = "'5.463"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.36%  "
$ws.Range("D20").Value2 = "'12.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.81%  "
$ws.Range("D21").Value = "2.100.21"
$ws.Range("E21").Value = "  -2.68%  "
$ws.Range("D22").Value2 = "'0.000007225"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.43%  "
$ws.Range("D23").Value2 = "'1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Value2 = "'6.126"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.89%  "
$ws.Range("D25").Value2 = "'167.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").Value2 = "'9.275"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.89%  "
$ws.Range("D27").Value2 = "'18.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.25%  "
$ws.Range("D28").Value2 = "'1.915"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.05%  "
$ws.Range("D29").Value2 = "'1.345"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.73%  "
$ws.Range("D30").Value2 = "'0.09587"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.89%  "
$ws.Range("D31").Value2 = "'4.406"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.72%  "
$ws.Range("E32").Value = "  -3.19%  "
$ws.Range("D33").Value2 = "'4.089"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.93%  "
$ws.Range("D34").Value2 = "'0.04664"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.53%  "
$ws.Range("E35").Value = "  -1.98%  "
$ws.Range("D36").Value2 = "'0.6943"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.48%  "
$ws.Range("D38").Value2 = "'2.709"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.48%  "
$ws.Range("D39").Value2 = "'0.01850"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.12%  "
$ws.Range("D40").Value2 = "'6.272"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.81%  "
$ws.Range("D41").Value2 = "'2.509"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.28%  "
$ws.Range("D42").Value2 = "'71.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.08%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value2 = "'1.943"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.87%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value2 = "'0.8563"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value2 = "'104.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.35%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value2 = "'1.000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").Value2 = "'0.4148"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.73%  "
$ws.Range("D48").Value = "1.021.93"
$ws.Range("E48").Value = "  +7.70%  "
$ws.Range("D49").Value2 = "'7.184"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.88%  "
$ws.Range("D50").Value2 = "'8.907"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.37%  "
$ws.Range("D51").Value2 = "'33.68"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.82%  "
